$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 199
$ws.Cells.Item(3, 6).Value = 674
$ws.Cells.Item(4, 6).Value = 673
$ws.Cells.Item(5, 6).Value = 533
$ws.Cells.Item(6, 6).Value = 2200
$ws.Cells.Item(7, 6).Value = 1297
$ws.Cells.Item(9, 6).Value = 80
$ws.Cells.Item(11, 6).Value = 2778
$ws.Cells.Item(17, 6).Value = 850
$ws.Cells.Item(19, 6).Value = 70
$ws.Cells.Item(21, 6).Value = 88
$ws.Cells.Item(22, 6).Value = 596
$ws.Cells.Item(26, 6).Value = 948
$ws.Cells.Item(27, 6).Value = 4839
$ws.Cells.Item(28, 6).Value = 347
$ws.Cells.Item(29, 6).Value = 112

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(21, 6).Value = 20
$ws.Cells.Item(25, 6).Value = 332
$ws.Cells.Item(27, 6).Value = 504
$ws.Cells.Item(31, 6).Value = 45
$ws.Cells.Item(37, 6).Value = 674

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 199
$ws.Cells.Item(7, 6).Value = 674
$ws.Cells.Item(9, 6).Value = 673
$ws.Cells.Item(12, 6).Value = 533
$ws.Cells.Item(13, 6).Value = 2200
$ws.Cells.Item(14, 6).Value = 1297
$ws.Cells.Item(16, 6).Value = 80
$ws.Cells.Item(20, 6).Value = 2778
$ws.Cells.Item(21, 6).Value = 2778
$ws.Cells.Item(30, 6).Value = 850
$ws.Cells.Item(31, 6).Value = 850
$ws.Cells.Item(33, 6).Value = 20
$ws.Cells.Item(35, 6).Value = 70
$ws.Cells.Item(36, 6).Value = 88
$ws.Cells.Item(39, 6).Value = 596
$ws.Cells.Item(41, 6).Value = 332
$ws.Cells.Item(42, 6).Value = 504
$ws.Cells.Item(45, 6).Value = 948
$ws.Cells.Item(46, 6).Value = 4839
$ws.Cells.Item(47, 6).Value = 45
$ws.Cells.Item(48, 6).Value = 347
$ws.Cells.Item(50, 6).Value = 112
$ws.Cells.Item(51, 6).Value = 674
